$d = $word.ActiveDocument
$W_NS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1. Remove the hidden "_GoBack" bookmark from the (currently empty) second
#    paragraph -- it reappears at the very end of the document (step 2e).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2. Append the new paragraphs after "Scipy is for Scientific Computations".

# 2a. "PyQT for cross-platform applications" (PyQT flagged by the spell
#     checker, same treatment as TensorFlow/Scipy above it)
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
[void]$target.InsertXML("<w:p xmlns:w='$W_NS'><w:proofErr w:type='spellStart'/><w:r><w:t>PyQT</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> for cross-platform applications</w:t></w:r></w:p>")

# 2b. empty paragraph
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
[void]$target.InsertXML("<w:p xmlns:w='$W_NS'></w:p>")

# 2c. "Pep8 – usually already integrated into the editor"
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$dash = [char]0x2013
[void]$target.InsertXML("<w:p xmlns:w='$W_NS'><w:r><w:t>Pep8 $dash usually already integrated into the editor</w:t></w:r></w:p>")

# 2d. "There needs to be 2 empty lines before a new class declaration"
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
[void]$target.InsertXML("<w:p xmlns:w='$W_NS'><w:r><w:t>There needs to be 2 empty lines before a new class declaration</w:t></w:r></w:p>")

# 2e. final empty paragraph holding the relocated "_GoBack" bookmark
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
[void]$target.InsertXML("<w:p xmlns:w='$W_NS'><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")
